# "Colocando header nos gráficos" - add a header label to column A of each
# data sheet (used as the category column for the charts), fix a handful of
# missing-accent labels, drop the unused "Teto" row on the emissions sheet,
# and refresh the cost-sheet header/values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheets 1-4 share the same row layout (Hidro..GD in A2:A12). Give each a
# header in A1 and fix the accented labels, then drop the old per-row
# bold/border style now that it lives on the header cell instead.
# ---------------------------------------------------------------------
$dataSheets = @(1, 2, 3, 4)

foreach ($idx in $dataSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # New header cell, formatted like the other row-1 header cells.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Fix accented labels.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # The label column no longer carries the bold/bordered style - that
    # moved up to the new header row.
    $ws.Range("A2:A12").Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet 5 - Emissoes Totais: header + accent fixes, and the "Teto" row
# (row 4) is removed entirely.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A2:A3").Style = "Normal"

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6 - Custo Total: header + accent fixes + updated figures.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value = "Tipo Expansão"

# B1's new label ("2015") must stay a text value like the other sheets'
# year headers, not get auto-converted to a number - paste the text value
# from an existing "2015" header cell instead of assigning a literal.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4163)

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 620
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
$ws6.Range("A2:A3").Style = "Normal"
